$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Maisy', 360, 'Active', 121, 13, 9, 3, 12, 440, 3.64, 18.04, '02. Straight Flush', '6,7,8,9,10,(D)', 15.09),
    @('Mark', 361, 'Active', 141, 15, 9, 0, 9, 507, 3.6, 14.12, '03. 4 of a Kind', 'Ace', 14.02),
    @('Matt', 362, 'Active', 189, 19, 11, 0, 11, 713, 3.77, 13.12, '03. 4 of a Kind', 'King', 16.02),
    @('Pepe', 364, 'Active', 103, 12, 9, 1, 10, 377, 3.66, 17.02, '03. 4 of a Kind', 'Ace', 15.08),
    @('Prashant', 365, 'Active', 36, 5, 9, 0, 9, 131, 3.64, 16.01, '03. 4 of a Kind', 'King', 16.01),
    @('Richard', 366, 'Active', 140, 19, 12, 0, 12, 605, 4.32, 12.05, '03. 4 of a Kind', 5, 18.05),
    @('Jon', 357, 'Active', 194, 19, 13, 0, 13, 697, 3.59, 13.08, '03. 4 of a Kind', 'Ace', 12.04),
    @('Alex', 348, 'Active', 83, 14, 9, 0, 9, 278, 3.35, 8.07, '03. 4 of a Kind', 'King', 11.02),
    @('Andy', 349, 'Active', 197, 19, 11, 0, 11, 817, 4.15, 11.04, '02. Straight Flush', '3,4,5,6,7', 19.07),
    @('Anthony', 350, 'Active', 119, 13, 10, 0, 10, 475, 3.99, 19.03, '02. Straight Flush', '7,8,9,10,Jack,(H)', 12.12),
    @('Bob', 351, 'Inactive', 35, 9, 11, 0, 11, 139, 3.97, 7.01, '04. Full House (3,2)', '4,5', 10.12),
    @('Chris', 353, 'Inactive', 1, 1, 3, 0, 3, 3, 3, 5.05, $null, $null, $null),
    @('Crafty', 354, 'Inactive', 17, 4, 11, 0, 11, 56, 3.29, 14.01, '03. 4 of a Kind', 2, 14.03),
    @('Illya', 355, 'Inactive', 28, 6, 10, 0, 10, 111, 3.96, 4.09, '05. Flush', 'Ace,4,8,9,King,(H)', 17.07),
    @('Jim', 356, 'Inactive', 3, 3, 9, 0, 9, 17, 5.67, 7.05, $null, $null, $null),
    @('John', 113, 'Inactive', 9, 3, 12, 0, 12, 46, 5.11, 5.06, $null, $null, $null),
    @('Jonathan', 358, 'Inactive', 10, 3, 7, 0, 7, 35, 3.5, 4.07, $null, $null, $null),
    @('Keith', 359, 'Inactive', 98, 12, 10, 0, 10, 333, 3.4, 13.04, '03. 4 of a Kind', 8, 12.06),
    @('Richard Snr', 367, 'Inactive', 1, 1, 1, 0, 1, 1, 1, 3.11, $null, $null, $null),
    @('Stuart', 368, 'Inactive', 48, 7, 12, 0, 12, 183, 3.81, 4.02, $null, $null, $null),
    @('Matthew', 363, 'Inactive', 2, 1, 3, 0, 3, 5, 2.5, 4.06, $null, $null, $null),
    @('Douglas', 424, 'Inactive', 5, 1, 7, 3, 10, 27, 5.4, 18.07, '04. Full House (3,2)', 'Ace,8', 18.07)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $rowData = $data[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $val = $rowData[$c - 1]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

Write-Host "Updated poker stats table (rows 2-23)"
